$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "thinBasic_Test_83789_912"

# Set column A width to match column D's width (25.7109375 in stored OOXML units).
# Excel COM's ColumnWidth uses "characters" (rounded to display pixels), so we use
# the character width that round-trips to the same stored width as column D.
$ws.Columns.Item(1).ColumnWidth = 24.8

# New "random" values for column B (B1:B20)
$newValues = @(1081, 1860, 1906, 34, 660, 1382, 588, 1272, 604, 1809, 575, 1947, 1276, 1431, 1124, 986, 72, 615, 788, 192)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

$excel.Calculate()
